# chore: update Sheets via scheduled runner
# Refreshes market-price-derived columns (currentAveragePrice*, LevePrice*,
# LeveProfit*) for a handful of leve rows across the ALC/ARM/BSM/CRP/CUL/GSM/
# LTW/WVR sheets, mirroring the data a scheduled price-refresh job would push.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 90.933334
$ws.Range("I2").Value = 87.63636
$ws.Range("K2").Value = 87.63636
$ws.Range("M2").Value = 25.36364

$ws.Range("H74").Value = 3945.4546
$ws.Range("J74").Value = 3983.3333
$ws.Range("L74").Value = 3983.3333
$ws.Range("N74").Value = -5855.3333

$ws.Range("H77").Value = 3945.4546
$ws.Range("J77").Value = 3983.3333
$ws.Range("L77").Value = 19916.6665
$ws.Range("N77").Value = -29276.6665

$ws.Range("H92").Value = 1284.1578
$ws.Range("I92").Value = 537.53845
$ws.Range("J92").Value = 2901.8333
$ws.Range("K92").Value = 537.53845
$ws.Range("L92").Value = 2901.8333
$ws.Range("M92").Value = 710.46155
$ws.Range("N92").Value = -5397.8333

$ws.Range("H135").Value = 22740.87
$ws.Range("I135").Value = 24576.38
$ws.Range("J135").Value = 3468
$ws.Range("K135").Value = 221187.42
$ws.Range("L135").Value = 31212
$ws.Range("M135").Value = -218652.42
$ws.Range("N135").Value = -36282

$ws.Range("H138").Value = 2429.2327
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 2429.2327
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 7287.6981
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = -17567.6981

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2003.3182
$ws.Range("I2").Value = 1592.8182
$ws.Range("J2").Value = 2413.818
$ws.Range("K2").Value = 1592.8182
$ws.Range("L2").Value = 2413.818
$ws.Range("M2").Value = -1479.8182
$ws.Range("N2").Value = -2639.818

$ws.Range("H13").Value = 23201500
$ws.Range("I13").Value = 29000000
$ws.Range("J13").Value = 7500
$ws.Range("K13").Value = 29000000
$ws.Range("L13").Value = 7500
$ws.Range("M13").Value = -28999856
$ws.Range("N13").Value = -7788

$ws.Range("H45").Value = 2550.6667
$ws.Range("I45").Value = 2754.1333
$ws.Range("J45").Value = 1533.3334
$ws.Range("K45").Value = 2754.1333
$ws.Range("L45").Value = 1533.3334
$ws.Range("M45").Value = -2377.1333
$ws.Range("N45").Value = -2287.3334

$ws.Range("H113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()

$ws.Range("H115").Value = 29583.334
$ws.Range("I115").Value = 25000
$ws.Range("K115").Value = 25000
$ws.Range("M115").Value = -23433

$ws.Range("H116").Value = 2003.3182
$ws.Range("I116").Value = 1592.8182
$ws.Range("J116").Value = 2413.818
$ws.Range("K116").Value = 1592.8182
$ws.Range("L116").Value = 2413.818
$ws.Range("M116").Value = 701.1818000000001
$ws.Range("N116").Value = -7001.818

$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()

$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()

$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()

$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()

$ws.Range("H121").Value = 65000
$ws.Range("J121").Value = 65000
$ws.Range("L121").Value = 65000
$ws.Range("N121").Value = -68494

$ws.Range("H122").Value = 7408475.5
$ws.Range("I122").Value = 1183.7273
$ws.Range("K122").Value = 3551.1819
$ws.Range("M122").Value = -1101.1819

$ws.Range("H132").Value = 176884.25
$ws.Range("I132").Value = 251324.75
$ws.Range("K132").Value = 753974.25
$ws.Range("M132").Value = -751444.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2003.3182
$ws.Range("I3").Value = 1592.8182
$ws.Range("J3").Value = 2413.818
$ws.Range("K3").Value = 1592.8182
$ws.Range("L3").Value = 2413.818
$ws.Range("M3").Value = -1478.8182
$ws.Range("N3").Value = -2641.818

$ws.Range("H64").Value = 891.6
$ws.Range("I64").Value = 503
$ws.Range("J64").Value = 988.75
$ws.Range("K64").Value = 503
$ws.Range("L64").Value = 988.75
$ws.Range("M64").Value = -278
$ws.Range("N64").Value = -1438.75

$ws.Range("H67").Value = 891.6
$ws.Range("I67").Value = 503
$ws.Range("J67").Value = 988.75
$ws.Range("K67").Value = 503
$ws.Range("L67").Value = 988.75
$ws.Range("M67").Value = 277
$ws.Range("N67").Value = -2548.75

$ws.Range("H94").Value = 674.7692
$ws.Range("I94").Value = 689.1111
$ws.Range("J94").Value = 642.5
$ws.Range("K94").Value = 689.1111
$ws.Range("L94").Value = 642.5
$ws.Range("M94").Value = -238.1111
$ws.Range("N94").Value = -1544.5

$ws.Range("H116").Value = 40000
$ws.Range("J116").Value = 40000
$ws.Range("L116").Value = 40000
$ws.Range("N116").Value = -49178

$ws.Range("H117").Value = 40000
$ws.Range("J117").Value = 40000
$ws.Range("L117").Value = 40000
$ws.Range("N117").Value = -49178

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 377.5
$ws.Range("J11").Value = 600
$ws.Range("L11").Value = 600
$ws.Range("N11").Value = -880

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 1390.4
$ws.Range("I11").Value = 1150
$ws.Range("J11").Value = 1751
$ws.Range("K11").Value = 3450
$ws.Range("L11").Value = 5253
$ws.Range("M11").Value = -3310
$ws.Range("N11").Value = -5533

$ws.Range("H102").Value = 3558.2856
$ws.Range("I102").Value = 2920
$ws.Range("J102").Value = 3912.889
$ws.Range("K102").Value = 8760
$ws.Range("L102").Value = 11738.667
$ws.Range("M102").Value = -6326
$ws.Range("N102").Value = -16606.667

$ws.Range("H113").Value = 573.5517
$ws.Range("I113").Value = 512.5
$ws.Range("J113").Value = 596.8095
$ws.Range("K113").Value = 1537.5
$ws.Range("L113").Value = 1790.4285
$ws.Range("M113").Value = 632.5
$ws.Range("N113").Value = -6130.4285

$ws.Range("H131").Value = 912.88464
$ws.Range("I131").Value = 498.3846
$ws.Range("J131").Value = 1051.0513
$ws.Range("K131").Value = 1495.1538
$ws.Range("L131").Value = 3153.1539
$ws.Range("M131").Value = 3544.8462
$ws.Range("N131").Value = -13233.1539

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H109").Value = 9785
$ws.Range("J109").Value = 9785
$ws.Range("L109").Value = 9785
$ws.Range("N109").Value = -11865

$ws.Range("H122").Value = 1845.6364
$ws.Range("I122").Value = 1911.75
$ws.Range("J122").Value = 1669.3334
$ws.Range("K122").Value = 5735.25
$ws.Range("L122").Value = 5008.0002
$ws.Range("M122").Value = -3285.25
$ws.Range("N122").Value = -9908.0002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 980
$ws.Range("I46").Value = 980
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 980
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -792
$ws.Range("N46").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 63865.03
$ws.Range("I132").Value = 42590.957
$ws.Range("J132").Value = 127687.25
$ws.Range("K132").Value = 127772.871
$ws.Range("L132").Value = 383061.75
$ws.Range("M132").Value = -125242.871
$ws.Range("N132").Value = -388121.75
